# Rename the worksheet from "Property1" to "DataNode" and move the
# active cell selection from AH12 to W37, matching the author's commit
# ("unify the conception of DataNode, DataTable, Entity.").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "DataNode"

$ws.Range("W37").Select()

$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 21
